$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.183885216712952
$ws.Range("B1").Value = 2.356223344802856
$ws.Range("C1").Value = 6.430662631988525
$ws.Range("D1").Value = 2.051105976104736
$ws.Range("E1").Value = 1.193235278129578
